$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46069
$ws.Range("B2").Value = 5520
$ws.Range("A3").Value = 46069.01041666666
$ws.Range("B3").Value = 5480
$ws.Range("A4").Value = 46069.02083333334
$ws.Range("B4").Value = 5450
$ws.Range("A5").Value = 46069.03125
$ws.Range("B5").Value = 5420
$ws.Range("A6").Value = 46069.04166666666
$ws.Range("B6").Value = 5420
$ws.Range("A7").Value = 46069.05208333334
$ws.Range("B7").Value = 5410
$ws.Range("A8").Value = 46069.0625
$ws.Range("B8").Value = 5400
$ws.Range("A9").Value = 46069.07291666666
$ws.Range("B9").Value = 5380
$ws.Range("A10").Value = 46069.08333333334
$ws.Range("B10").Value = 5370
$ws.Range("A11").Value = 46069.09375
$ws.Range("B11").Value = 5360
$ws.Range("A12").Value = 46069.10416666666
$ws.Range("B12").Value = 5360
$ws.Range("A13").Value = 46069.11458333334
$ws.Range("B13").Value = 5370
$ws.Range("A14").Value = 46069.125
$ws.Range("B14").Value = 5380
$ws.Range("A15").Value = 46069.13541666666
$ws.Range("B15").Value = 5400
$ws.Range("A16").Value = 46069.14583333334
$ws.Range("B16").Value = 5420
$ws.Range("A17").Value = 46069.15625
$ws.Range("B17").Value = 5480
$ws.Range("A18").Value = 46069.16666666666
$ws.Range("B18").Value = 5570
$ws.Range("A19").Value = 46069.17708333334
$ws.Range("B19").Value = 5680
$ws.Range("A20").Value = 46069.1875
$ws.Range("B20").Value = 5810
$ws.Range("A21").Value = 46069.19791666666
$ws.Range("B21").Value = 5960
$ws.Range("A22").Value = 46069.20833333334
$ws.Range("B22").Value = 6140
$ws.Range("A23").Value = 46069.21875
$ws.Range("B23").Value = 6330
$ws.Range("A24").Value = 46069.22916666666
$ws.Range("B24").Value = 6540
$ws.Range("A25").Value = 46069.23958333334
$ws.Range("B25").Value = 6760
$ws.Range("A26").Value = 46069.25
$ws.Range("B26").Value = 7000
$ws.Range("A27").Value = 46069.26041666666
$ws.Range("B27").Value = 7220
$ws.Range("A28").Value = 46069.27083333334
$ws.Range("B28").Value = 7430
$ws.Range("A29").Value = 46069.28125
$ws.Range("B29").Value = 7650
$ws.Range("A30").Value = 46069.29166666666
$ws.Range("B30").Value = 7850
$ws.Range("A31").Value = 46069.30208333334
$ws.Range("B31").Value = 8010
$ws.Range("A32").Value = 46069.3125
$ws.Range("B32").Value = 8100
$ws.Range("A33").Value = 46069.32291666666
$ws.Range("B33").Value = 8160
$ws.Range("A34").Value = 46069.33333333334
$ws.Range("B34").Value = 8180
$ws.Range("A35").Value = 46069.34375
$ws.Range("B35").Value = 8170
$ws.Range("A36").Value = 46069.35416666666
$ws.Range("B36").Value = 8140
$ws.Range("A37").Value = 46069.36458333334
$ws.Range("B37").Value = 8070
$ws.Range("A38").Value = 46069.375
$ws.Range("B38").Value = 7980
$ws.Range("A39").Value = 46069.38541666666
$ws.Range("B39").Value = 7880
$ws.Range("A40").Value = 46069.39583333334
$ws.Range("B40").Value = 7780
$ws.Range("A41").Value = 46069.40625
$ws.Range("B41").Value = 7690
$ws.Range("A42").Value = 46069.41666666666
$ws.Range("B42").Value = 7600
$ws.Range("A43").Value = 46069.42708333334
$ws.Range("B43").Value = 7530
$ws.Range("A44").Value = 46069.4375
$ws.Range("B44").Value = 7470
$ws.Range("A45").Value = 46069.44791666666
$ws.Range("B45").Value = 7420
$ws.Range("A46").Value = 46069.45833333334
$ws.Range("B46").Value = 7350
$ws.Range("A47").Value = 46069.46875
$ws.Range("B47").Value = 7310
$ws.Range("A48").Value = 46069.47916666666
$ws.Range("B48").Value = 7290
$ws.Range("A49").Value = 46069.48958333334
$ws.Range("B49").Value = 7280
$ws.Range("A50").Value = 46069.5
$ws.Range("B50").Value = 7280
$ws.Range("A51").Value = 46069.51041666666
$ws.Range("B51").Value = 7270
$ws.Range("A52").Value = 46069.52083333334
$ws.Range("B52").Value = 7270
$ws.Range("A53").Value = 46069.53125
$ws.Range("B53").Value = 7270
$ws.Range("A54").Value = 46069.54166666666
$ws.Range("B54").Value = 7270
$ws.Range("A55").Value = 46069.55208333334
$ws.Range("B55").Value = 7270
$ws.Range("A56").Value = 46069.5625
$ws.Range("B56").Value = 7280
$ws.Range("A57").Value = 46069.57291666666
$ws.Range("B57").Value = 7310
$ws.Range("A58").Value = 46069.58333333334
$ws.Range("B58").Value = 7340
$ws.Range("A59").Value = 46069.59375
$ws.Range("B59").Value = 7380
$ws.Range("A60").Value = 46069.60416666666
$ws.Range("B60").Value = 7410
$ws.Range("A61").Value = 46069.61458333334
$ws.Range("B61").Value = 7460
$ws.Range("A62").Value = 46069.625
$ws.Range("B62").Value = 7540
$ws.Range("A63").Value = 46069.63541666666
$ws.Range("B63").Value = 7610
$ws.Range("A64").Value = 46069.64583333334
$ws.Range("B64").Value = 7700
$ws.Range("A65").Value = 46069.65625
$ws.Range("B65").Value = 7780
$ws.Range("A66").Value = 46069.66666666666
$ws.Range("B66").Value = 7880
$ws.Range("A67").Value = 46069.67708333334
$ws.Range("B67").Value = 7990
$ws.Range("A68").Value = 46069.6875
$ws.Range("B68").Value = 8080
$ws.Range("A69").Value = 46069.69791666666
$ws.Range("B69").Value = 8170
$ws.Range("A70").Value = 46069.70833333334
$ws.Range("B70").Value = 8270
$ws.Range("A71").Value = 46069.71875
$ws.Range("B71").Value = 8320
$ws.Range("A72").Value = 46069.72916666666
$ws.Range("B72").Value = 8340
$ws.Range("A73").Value = 46069.73958333334
$ws.Range("B73").Value = 8360
$ws.Range("A74").Value = 46069.75
$ws.Range("B74").Value = 8360
$ws.Range("A75").Value = 46069.76041666666
$ws.Range("B75").Value = 8340
$ws.Range("A76").Value = 46069.77083333334
$ws.Range("B76").Value = 8290
$ws.Range("A77").Value = 46069.78125
$ws.Range("B77").Value = 8240
$ws.Range("A78").Value = 46069.79166666666
$ws.Range("B78").Value = 8180
$ws.Range("A79").Value = 46069.80208333334
$ws.Range("B79").Value = 8130
$ws.Range("A80").Value = 46069.8125
$ws.Range("B80").Value = 8070
$ws.Range("A81").Value = 46069.82291666666
$ws.Range("B81").Value = 7990
$ws.Range("A82").Value = 46069.83333333334
$ws.Range("B82").Value = 7860
$ws.Range("A83").Value = 46069.84375
$ws.Range("B83").Value = 7730
$ws.Range("A84").Value = 46069.85416666666
$ws.Range("B84").Value = 7630
$ws.Range("A85").Value = 46069.86458333334
$ws.Range("B85").Value = 7480
$ws.Range("A86").Value = 46069.875
$ws.Range("B86").Value = 7300
$ws.Range("A87").Value = 46069.88541666666
$ws.Range("B87").Value = 7160
$ws.Range("A88").Value = 46069.89583333334
$ws.Range("B88").Value = 7010
$ws.Range("A89").Value = 46069.90625
$ws.Range("B89").Value = 6880
$ws.Range("A90").Value = 46069.91666666666
$ws.Range("B90").Value = 6720
$ws.Range("A91").Value = 46069.92708333334
$ws.Range("B91").Value = 6600
$ws.Range("A92").Value = 46069.9375
$ws.Range("B92").Value = 6500
$ws.Range("A93").Value = 46069.94791666666
$ws.Range("B93").Value = 6380
$ws.Range("A94").Value = 46069.95833333334
$ws.Range("B94").Value = 6250
$ws.Range("A95").Value = 46069.96875
$ws.Range("B95").Value = 6140
$ws.Range("A96").Value = 46069.97916666666
$ws.Range("B96").Value = 6060
$ws.Range("A97").Value = 46069.98958333334
$ws.Range("B97").Value = 6010
